$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2049
$ws.Range("J17").Value = 2049
$ws.Range("L17").Value = 6147
$ws.Range("N17").Value = -6483
$ws.Range("H18").Value = 2031.9286
$ws.Range("I18").Value = 1270.875
$ws.Range("J18").Value = 3046.6667
$ws.Range("K18").Value = 1270.875
$ws.Range("L18").Value = 3046.6667
$ws.Range("M18").Value = -986.875
$ws.Range("N18").Value = -3614.6667
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H33").Value = 2709.4614
$ws.Range("I33").Value = 2927.5833
$ws.Range("K33").Value = 2927.5833
$ws.Range("M33").Value = -2698.5833
$ws.Range("H53").Value = 5077.6665
$ws.Range("I53").Value = 6782.125
$ws.Range("K53").Value = 6782.125
$ws.Range("M53").Value = -6145.125
$ws.Range("H59").Value = 3240.8572
$ws.Range("I59").Value = 900
$ws.Range("K59").Value = 2700
$ws.Range("M59").Value = -2143
$ws.Range("H69").Value = 775485.6
$ws.Range("I69").Value = 775485.6
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 2326456.8
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -2325582.8
$ws.Range("H70").Value = 1932.3334
$ws.Range("J70").Value = 2917.875
$ws.Range("L70").Value = 8753.625
$ws.Range("N70").Value = -9293.625
$ws.Range("H72").Value = 775485.6
$ws.Range("I72").Value = 775485.6
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 6979370.399999999
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -6975002.399999999
$ws.Range("H73").Value = 1932.3334
$ws.Range("J73").Value = 2917.875
$ws.Range("L73").Value = 8753.625
$ws.Range("N73").Value = -10625.625
$ws.Range("H74").Value = 13437.5
$ws.Range("I74").Value = 6666.6665
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 6666.6665
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -5730.6665
$ws.Range("N74").Value = -16872
$ws.Range("H77").Value = 13437.5
$ws.Range("I77").Value = 6666.6665
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 33333.3325
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -28653.3325
$ws.Range("N77").Value = -84360
$ws.Range("H80").Value = 1020.2308
$ws.Range("I80").Value = 1077.6
$ws.Range("J80").Value = 829
$ws.Range("K80").Value = 3232.8
$ws.Range("L80").Value = 2487
$ws.Range("M80").Value = -2234.8
$ws.Range("N80").Value = -4483
$ws.Range("H83").Value = 1020.2308
$ws.Range("I83").Value = 1077.6
$ws.Range("J83").Value = 829
$ws.Range("K83").Value = 9698.4
$ws.Range("L83").Value = 7461
$ws.Range("M83").Value = -4706.4
$ws.Range("N83").Value = -17445
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
$ws.Range("H88").Value = 17584.297
$ws.Range("I88").Value = 2282.5715
$ws.Range("J88").Value = 22939.9
$ws.Range("K88").Value = 2282.5715
$ws.Range("L88").Value = 22939.9
$ws.Range("M88").Value = -1876.5715
$ws.Range("N88").Value = -23751.9
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384
$ws.Range("H91").Value = 17584.297
$ws.Range("I91").Value = 2282.5715
$ws.Range("J91").Value = 22939.9
$ws.Range("K91").Value = 2282.5715
$ws.Range("L91").Value = 22939.9
$ws.Range("M91").Value = -878.5715
$ws.Range("N91").Value = -25747.9
$ws.Range("H92").Value = 128945.82
$ws.Range("I92").Value = 56490
$ws.Range("J92").Value = 454997
$ws.Range("K92").Value = 56490
$ws.Range("L92").Value = 454997
$ws.Range("M92").Value = -55242
$ws.Range("N92").Value = -457493
$ws.Range("H94").Value = 1394
$ws.Range("I94").Value = 1394
$ws.Range("K94").Value = 1394
$ws.Range("M94").Value = -943
$ws.Range("H97").Value = 972.2308
$ws.Range("I97").Value = 853.5
$ws.Range("J97").Value = 993.8182
$ws.Range("K97").Value = 2560.5
$ws.Range("L97").Value = 2981.4546
$ws.Range("M97").Value = -2064.5
$ws.Range("N97").Value = -3973.4546
$ws.Range("H99").Value = 378.2857
$ws.Range("I99").Value = 374.66666
$ws.Range("K99").Value = 1123.99998
$ws.Range("M99").Value = 374.0000199999999
$ws.Range("I100").Value = 8888.333000000001
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 8888.333000000001
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -8347.333000000001
$ws.Range("H101").Value = 3079.3333
$ws.Range("I101").Value = 2707.2
$ws.Range("J101").Value = 3345.1428
$ws.Range("K101").Value = 8121.599999999999
$ws.Range("L101").Value = 10035.4284
$ws.Range("M101").Value = -6499.599999999999
$ws.Range("N101").Value = -13279.4284
$ws.Range("H103").Value = 3315.8462
$ws.Range("J103").Value = 3780.5454
$ws.Range("L103").Value = 11341.6362
$ws.Range("N103").Value = -12513.6362
$ws.Range("H112").Value = 38010.24
$ws.Range("J112").Value = 47217.95
$ws.Range("L112").Value = 141653.85
$ws.Range("N112").Value = -143869.85
$ws.Range("H115").Value = 963.2632
$ws.Range("I115").Value = 735.4706
$ws.Range("K115").Value = 2206.4118
$ws.Range("M115").Value = -639.4117999999999
$ws.Range("H129").Value = 27981.5
$ws.Range("I129").Value = 65586.64
$ws.Range("J129").Value = 4050.9546
$ws.Range("K129").Value = 196759.92
$ws.Range("L129").Value = 12152.8638
$ws.Range("M129").Value = -191759.92
$ws.Range("N129").Value = -22152.8638
$ws.Range("H137").Value = 2393.325
$ws.Range("I137").Value = 1630.7097
$ws.Range("J137").Value = 5020.1113
$ws.Range("K137").Value = 4892.1291
$ws.Range("L137").Value = 15060.3339
$ws.Range("M137").Value = -2342.1291
$ws.Range("N137").Value = -20160.3339
$ws.Range("H138").Value = 3266.4614
$ws.Range("I138").Value = 1127.1305
$ws.Range("J138").Value = 4438
$ws.Range("K138").Value = 3381.3915
$ws.Range("L138").Value = 13314
$ws.Range("M138").Value = 1758.6085
$ws.Range("N138").Value = -23594
$ws.Range("H141").Value = 3259.6
$ws.Range("I141").Value = 2099.6667
$ws.Range("J141").Value = 4999.5
$ws.Range("K141").Value = 6299.000100000001
$ws.Range("L141").Value = 14998.5
$ws.Range("M141").Value = -1119.000100000001
$ws.Range("N141").Value = -25358.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24520.178
$ws.Range("I32").Value = 29385.5
$ws.Range("J32").Value = 15674.137
$ws.Range("K32").Value = 29385.5
$ws.Range("L32").Value = 15674.137
$ws.Range("M32").Value = -29098.5
$ws.Range("N32").Value = -16248.137
$ws.Range("H45").Value = 1472.1666
$ws.Range("I45").Value = 1093.5
$ws.Range("K45").Value = 1093.5
$ws.Range("M45").Value = -716.5
$ws.Range("H64").Value = 49979.75
$ws.Range("J64").Value = 49979.75
$ws.Range("L64").Value = 49979.75
$ws.Range("N64").Value = -50475.75
$ws.Range("H67").Value = 49979.75
$ws.Range("J67").Value = 49979.75
$ws.Range("L67").Value = 49979.75
$ws.Range("N67").Value = -51695.75
$ws.Range("H74").Value = 204094.11
$ws.Range("I74").Value = 183395.4
$ws.Range("K74").Value = 183395.4
$ws.Range("M74").Value = -182521.4
$ws.Range("H77").Value = 204094.11
$ws.Range("I77").Value = 183395.4
$ws.Range("K77").Value = 916977
$ws.Range("M77").Value = -912609
$ws.Range("H97").Value = 1033
$ws.Range("I97").Value = 917.75
$ws.Range("J97").Value = 1494
$ws.Range("K97").Value = 917.75
$ws.Range("L97").Value = 1494
$ws.Range("M97").Value = -421.75
$ws.Range("N97").Value = -2486
$ws.Range("H102").Value = 1668499.6
$ws.Range("I102").Value = 2001599.6
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2001599.6
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1999977.6
$ws.Range("N102").Value = -6244
$ws.Range("H122").Value = 38899.32
$ws.Range("I122").Value = 2828.7646
$ws.Range("J122").Value = 115549.25
$ws.Range("K122").Value = 8486.293799999999
$ws.Range("L122").Value = 346647.75
$ws.Range("M122").Value = -6036.293799999999
$ws.Range("N122").Value = -351547.75
$ws.Range("H132").Value = 3450.9614
$ws.Range("I132").Value = 3162.3914
$ws.Range("J132").Value = 5663.3335
$ws.Range("K132").Value = 9487.174199999999
$ws.Range("L132").Value = 16990.0005
$ws.Range("M132").Value = -6957.174199999999
$ws.Range("N132").Value = -22050.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2448.5
$ws.Range("I8").Value = 2448.5
$ws.Range("K8").Value = 2448.5
$ws.Range("M8").Value = -2308.5
$ws.Range("H20").Value = 5391.778
$ws.Range("I20").Value = 4851.6
$ws.Range("J20").Value = 6067
$ws.Range("K20").Value = 4851.6
$ws.Range("L20").Value = 6067
$ws.Range("M20").Value = -4604.6
$ws.Range("N20").Value = -6561
$ws.Range("H60").Value = 49999
$ws.Range("J60").Value = 49999
$ws.Range("L60").Value = 49999
$ws.Range("N60").Value = -51197
$ws.Range("H86").Value = 62047.535
$ws.Range("I86").Value = 2713.625
$ws.Range("J86").Value = 129857.71
$ws.Range("K86").Value = 2713.625
$ws.Range("L86").Value = 129857.71
$ws.Range("M86").Value = -1590.625
$ws.Range("N86").Value = -132103.71
$ws.Range("H89").Value = 62047.535
$ws.Range("I89").Value = 2713.625
$ws.Range("J89").Value = 129857.71
$ws.Range("K89").Value = 13568.125
$ws.Range("L89").Value = 649288.55
$ws.Range("M89").Value = -7952.125
$ws.Range("N89").Value = -660520.55
$ws.Range("H99").Value = 2647.3
$ws.Range("I99").Value = 2113.1667
$ws.Range("J99").Value = 3448.5
$ws.Range("K99").Value = 2113.1667
$ws.Range("L99").Value = 3448.5
$ws.Range("M99").Value = -615.1667000000002
$ws.Range("N99").Value = -6444.5
$ws.Range("H105").Value = 2273.4
$ws.Range("I105").Value = 1499.5
$ws.Range("K105").Value = 1499.5
$ws.Range("M105").Value = 247.5
$ws.Range("H134").Value = 25182.59
$ws.Range("I134").Value = 30777.709
$ws.Range("J134").Value = 3501.5
$ws.Range("K134").Value = 92333.12699999999
$ws.Range("L134").Value = 10504.5
$ws.Range("M134").Value = -89798.12699999999
$ws.Range("N134").Value = -15574.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1917.6
$ws.Range("I16").Value = 2026
$ws.Range("J16").Value = 1664.6666
$ws.Range("K16").Value = 2026
$ws.Range("L16").Value = 1664.6666
$ws.Range("M16").Value = -1739
$ws.Range("N16").Value = -2238.6666
$ws.Range("H22").Value = 787.1
$ws.Range("I22").Value = 538.6
$ws.Range("J22").Value = 1035.6
$ws.Range("K22").Value = 538.6
$ws.Range("L22").Value = 1035.6
$ws.Range("M22").Value = -188.6
$ws.Range("N22").Value = -1735.6
$ws.Range("H31").Value = 2358.7256
$ws.Range("I31").Value = 1988.9143
$ws.Range("J31").Value = 3167.6875
$ws.Range("K31").Value = 1988.9143
$ws.Range("L31").Value = 3167.6875
$ws.Range("M31").Value = -1693.9143
$ws.Range("N31").Value = -3757.6875
$ws.Range("H34").Value = 2358.7256
$ws.Range("I34").Value = 1988.9143
$ws.Range("J34").Value = 3167.6875
$ws.Range("K34").Value = 1988.9143
$ws.Range("L34").Value = 3167.6875
$ws.Range("M34").Value = -1786.9143
$ws.Range("N34").Value = -3571.6875
$ws.Range("H58").Value = 9725.4
$ws.Range("I58").Value = 9499.5
$ws.Range("J58").Value = 9876
$ws.Range("K58").Value = 9499.5
$ws.Range("L58").Value = 9876
$ws.Range("M58").Value = -9296.5
$ws.Range("N58").Value = -10282
$ws.Range("H62").Value = 84621.69500000001
$ws.Range("I62").Value = 255751.25
$ws.Range("J62").Value = 8564.111000000001
$ws.Range("K62").Value = 255751.25
$ws.Range("L62").Value = 8564.111000000001
$ws.Range("M62").Value = -255127.25
$ws.Range("N62").Value = -9812.111000000001
$ws.Range("H65").Value = 84621.69500000001
$ws.Range("I65").Value = 255751.25
$ws.Range("J65").Value = 8564.111000000001
$ws.Range("K65").Value = 1278756.25
$ws.Range("L65").Value = 42820.55500000001
$ws.Range("M65").Value = -1275636.25
$ws.Range("N65").Value = -49060.55500000001
$ws.Range("H94").Value = 2263.8572
$ws.Range("I94").Value = 3580
$ws.Range("J94").Value = 509
$ws.Range("K94").Value = 3580
$ws.Range("L94").Value = 509
$ws.Range("M94").Value = -3129
$ws.Range("N94").Value = -1411
$ws.Range("H113").Value = 1917.6
$ws.Range("I113").Value = 2026
$ws.Range("J113").Value = 1664.6666
$ws.Range("K113").Value = 2026
$ws.Range("L113").Value = 1664.6666
$ws.Range("M113").Value = 144
$ws.Range("N113").Value = -6004.6666
$ws.Range("H132").Value = 2926.5557
$ws.Range("I132").Value = 3057.9412
$ws.Range("J132").Value = 693
$ws.Range("K132").Value = 9173.8236
$ws.Range("L132").Value = 2079
$ws.Range("M132").Value = -6643.8236
$ws.Range("N132").Value = -7139
$ws.Range("H134").Value = 5702.4287
$ws.Range("I134").Value = 5381.385
$ws.Range("K134").Value = 16144.155
$ws.Range("M134").Value = -13609.155
$ws.Range("H135").Value = 67966.664
$ws.Range("J135").Value = 67966.664
$ws.Range("L135").Value = 67966.664
$ws.Range("N135").Value = -78106.664
$ws.Range("H136").Value = 9725.4
$ws.Range("I136").Value = 9499.5
$ws.Range("J136").Value = 9876
$ws.Range("K136").Value = 28498.5
$ws.Range("L136").Value = 29628
$ws.Range("M136").Value = -25948.5
$ws.Range("N136").Value = -34728
$ws.Range("H141").Value = 189566.34
$ws.Range("J141").Value = 189566.34
$ws.Range("L141").Value = 189566.34
$ws.Range("N141").Value = -199926.34

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 481.57144
$ws.Range("I5").Value = 407.82352
$ws.Range("J5").Value = 795
$ws.Range("K5").Value = 1223.47056
$ws.Range("L5").Value = 2385
$ws.Range("M5").Value = -1111.47056
$ws.Range("N5").Value = -2609
$ws.Range("H14").Value = 1738.3043
$ws.Range("I14").Value = 1738.3043
$ws.Range("K14").Value = 5214.9129
$ws.Range("M14").Value = -5041.9129
$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 300
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 900
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -670
$ws.Range("H29").Value = 227987.6
$ws.Range("I29").Value = 534709.5
$ws.Range("J29").Value = 23506.334
$ws.Range("K29").Value = 1604128.5
$ws.Range("L29").Value = 70519.00199999999
$ws.Range("M29").Value = -1603851.5
$ws.Range("N29").Value = -71073.00199999999
$ws.Range("H46").Value = 1493.5
$ws.Range("I46").Value = 331.33334
$ws.Range("J46").Value = 1991.5714
$ws.Range("K46").Value = 994.0000200000001
$ws.Range("L46").Value = 5974.7142
$ws.Range("M46").Value = -903.0000200000001
$ws.Range("N46").Value = -6156.7142
$ws.Range("H51").Value = 2186.5
$ws.Range("I51").Value = 2229.75
$ws.Range("J51").Value = 2100
$ws.Range("K51").Value = 6689.25
$ws.Range("L51").Value = 6300
$ws.Range("M51").Value = -6229.25
$ws.Range("N51").Value = -7220
$ws.Range("H56").Value = 7299.375
$ws.Range("I56").Value = 7299.375
$ws.Range("K56").Value = 7299.375
$ws.Range("M56").Value = -6769.375
$ws.Range("H69").Value = 4940
$ws.Range("J69").Value = 4925
$ws.Range("L69").Value = 14775
$ws.Range("N69").Value = -16397
$ws.Range("H72").Value = 4940
$ws.Range("J72").Value = 4925
$ws.Range("L72").Value = 44325
$ws.Range("N72").Value = -52437
$ws.Range("H80").Value = 17558.8
$ws.Range("I80").Value = 7996
$ws.Range("J80").Value = 21657.143
$ws.Range("K80").Value = 23988
$ws.Range("L80").Value = 64971.429
$ws.Range("M80").Value = -23052
$ws.Range("N80").Value = -66843.429
$ws.Range("H83").Value = 17558.8
$ws.Range("I83").Value = 7996
$ws.Range("J83").Value = 21657.143
$ws.Range("K83").Value = 71964
$ws.Range("L83").Value = 194914.287
$ws.Range("M83").Value = -67284
$ws.Range("N83").Value = -204274.287
$ws.Range("H86").Value = 259.6
$ws.Range("I86").Value = 177.55556
$ws.Range("J86").Value = 998
$ws.Range("K86").Value = 532.66668
$ws.Range("L86").Value = 2994
$ws.Range("M86").Value = 653.33332
$ws.Range("N86").Value = -5366
$ws.Range("H89").Value = 259.6
$ws.Range("I89").Value = 177.55556
$ws.Range("J89").Value = 998
$ws.Range("K89").Value = 1598.00004
$ws.Range("L89").Value = 8982
$ws.Range("M89").Value = 4329.99996
$ws.Range("N89").Value = -20838
$ws.Range("H98").Value = 1874.4
$ws.Range("I98").Value = 1569
$ws.Range("J98").Value = 1950.75
$ws.Range("K98").Value = 4707
$ws.Range("L98").Value = 5852.25
$ws.Range("M98").Value = -3209
$ws.Range("N98").Value = -8848.25
$ws.Range("H107").Value = 1367.9678
$ws.Range("J107").Value = 1472.7858
$ws.Range("L107").Value = 4418.357400000001
$ws.Range("N107").Value = -8258.357400000001
$ws.Range("H131").Value = 3881.7334
$ws.Range("J131").Value = 5242.8887
$ws.Range("L131").Value = 15728.6661
$ws.Range("N131").Value = -25808.6661
$ws.Range("H135").Value = 481.57144
$ws.Range("I135").Value = 407.82352
$ws.Range("J135").Value = 795
$ws.Range("K135").Value = 3670.41168
$ws.Range("L135").Value = 7155
$ws.Range("M135").Value = -1135.41168
$ws.Range("N135").Value = -12225
$ws.Range("H137").Value = 109251.77
$ws.Range("I137").Value = 1434412.8
$ws.Range("J137").Value = 5025.6294
$ws.Range("K137").Value = 4303238.4
$ws.Range("L137").Value = 15076.8882
$ws.Range("M137").Value = -4298138.4
$ws.Range("N137").Value = -25276.8882
$ws.Range("H140").Value = 3386.4375
$ws.Range("I140").Value = 1550.0834
$ws.Range("J140").Value = 4488.25
$ws.Range("K140").Value = 4650.2502
$ws.Range("L140").Value = 13464.75
$ws.Range("M140").Value = 529.7497999999996
$ws.Range("N140").Value = -23824.75
$ws.Range("H141").Value = 3673.7856
$ws.Range("I141").Value = 3673.7856
$ws.Range("K141").Value = 11021.3568
$ws.Range("M141").Value = -5841.356800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 39800
$ws.Range("J15").Value = 39800
$ws.Range("L15").Value = 39800
$ws.Range("N15").Value = -40376
$ws.Range("H43").Value = 15242.5
$ws.Range("I43").Value = 4651.6665
$ws.Range("K43").Value = 4651.6665
$ws.Range("M43").Value = -4500.6665
$ws.Range("H80").Value = 22715.834
$ws.Range("I80").Value = 1347.5
$ws.Range("J80").Value = 33400
$ws.Range("K80").Value = 1347.5
$ws.Range("L80").Value = 33400
$ws.Range("M80").Value = -349.5
$ws.Range("N80").Value = -35396
$ws.Range("H81").Value = 39800
$ws.Range("J81").Value = 39800
$ws.Range("L81").Value = 39800
$ws.Range("N81").Value = -41796
$ws.Range("H83").Value = 22715.834
$ws.Range("I83").Value = 1347.5
$ws.Range("J83").Value = 33400
$ws.Range("K83").Value = 6737.5
$ws.Range("L83").Value = 167000
$ws.Range("M83").Value = -1745.5
$ws.Range("N83").Value = -176984
$ws.Range("H84").Value = 39800
$ws.Range("J84").Value = 39800
$ws.Range("L84").Value = 119400
$ws.Range("N84").Value = -129384
$ws.Range("H97").Value = 901.0833
$ws.Range("I97").Value = 870.55554
$ws.Range("J97").Value = 992.6667
$ws.Range("K97").Value = 870.55554
$ws.Range("L97").Value = 992.6667
$ws.Range("M97").Value = -374.55554
$ws.Range("N97").Value = -1984.6667
$ws.Range("H122").Value = 3467.3076
$ws.Range("I122").Value = 2795.0833
$ws.Range("J122").Value = 4043.5
$ws.Range("K122").Value = 8385.249899999999
$ws.Range("L122").Value = 12130.5
$ws.Range("M122").Value = -5935.249899999999
$ws.Range("N122").Value = -17030.5
$ws.Range("H126").Value = 65084.656
$ws.Range("I126").Value = 89681.60000000001
$ws.Range("K126").Value = 269044.8
$ws.Range("M126").Value = -266574.8
$ws.Range("H132").Value = 837272.7
$ws.Range("I132").Value = 979540.2
$ws.Range("J132").Value = 3991.7144
$ws.Range("K132").Value = 2938620.6
$ws.Range("L132").Value = 11975.1432
$ws.Range("M132").Value = -2936090.6
$ws.Range("N132").Value = -17035.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 24999.143
$ws.Range("I3").Value = 24999
$ws.Range("J3").Value = 24999.25
$ws.Range("K3").Value = 24999
$ws.Range("L3").Value = 24999.25
$ws.Range("M3").Value = -24887
$ws.Range("N3").Value = -25223.25
$ws.Range("H9").Value = 2479.375
$ws.Range("I9").Value = 1973.1666
$ws.Range("J9").Value = 3998
$ws.Range("K9").Value = 1973.1666
$ws.Range("L9").Value = 3998
$ws.Range("M9").Value = -1749.1666
$ws.Range("N9").Value = -4446
$ws.Range("H11").Value = 3500
$ws.Range("J11").Value = 3500
$ws.Range("L11").Value = 3500
$ws.Range("N11").Value = -3780
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H13").Value = 1780.25
$ws.Range("I13").Value = 2226.3635
$ws.Range("J13").Value = 798.8
$ws.Range("K13").Value = 2226.3635
$ws.Range("L13").Value = 798.8
$ws.Range("M13").Value = -2086.3635
$ws.Range("N13").Value = -1078.8
$ws.Range("H14").Value = 24000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 24000
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").Value = 24000
$ws.Range("N14").Value = -24344
$ws.Range("H15").Value = 24999.143
$ws.Range("I15").Value = 24999
$ws.Range("J15").Value = 24999.25
$ws.Range("K15").Value = 24999
$ws.Range("L15").Value = 24999.25
$ws.Range("M15").Value = -24829
$ws.Range("N15").Value = -25339.25
$ws.Range("H16").Value = 2797.818
$ws.Range("I16").Value = 1475.1111
$ws.Range("K16").Value = 1475.1111
$ws.Range("M16").Value = -1305.1111
$ws.Range("H20").Value = 9272.056
$ws.Range("I20").Value = 7142.8096
$ws.Range("J20").Value = 12253
$ws.Range("K20").Value = 7142.8096
$ws.Range("L20").Value = 12253
$ws.Range("M20").Value = -6916.8096
$ws.Range("N20").Value = -12705
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = 0
$ws.Range("H25").Value = 18599
$ws.Range("J25").Value = 18599
$ws.Range("L25").Value = 18599
$ws.Range("N25").Value = -19059
$ws.Range("H40").Value = 40201.246
$ws.Range("I40").Value = 46144.746
$ws.Range("K40").Value = 46144.746
$ws.Range("M40").Value = -46008.746
$ws.Range("H61").Value = 4004.8948
$ws.Range("I61").Value = 3240.4707
$ws.Range("J61").Value = 10502.5
$ws.Range("K61").Value = 3240.4707
$ws.Range("L61").Value = 10502.5
$ws.Range("M61").Value = -3038.4707
$ws.Range("N61").Value = -10906.5
$ws.Range("H70").Value = 52498
$ws.Range("J70").Value = 52498
$ws.Range("L70").Value = 52498
$ws.Range("N70").Value = -53038
$ws.Range("H73").Value = 52498
$ws.Range("J73").Value = 52498
$ws.Range("L73").Value = 52498
$ws.Range("N73").Value = -54370
$ws.Range("H74").Value = 79842.60000000001
$ws.Range("I74").Value = 73332.336
$ws.Range("J74").Value = 89608
$ws.Range("K74").Value = 73332.336
$ws.Range("L74").Value = 89608
$ws.Range("M74").Value = -72334.336
$ws.Range("N74").Value = -91604
$ws.Range("H77").Value = 79842.60000000001
$ws.Range("I77").Value = 73332.336
$ws.Range("J77").Value = 89608
$ws.Range("K77").Value = 219997.008
$ws.Range("L77").Value = 268824
$ws.Range("M77").Value = -215005.008
$ws.Range("N77").Value = -278808
$ws.Range("H80").Value = 65888.664
$ws.Range("I80").Value = 56000
$ws.Range("J80").Value = 78249.5
$ws.Range("K80").Value = 56000
$ws.Range("L80").Value = 78249.5
$ws.Range("M80").Value = -54877
$ws.Range("N80").Value = -80495.5
$ws.Range("H82").Value = 3500.4285
$ws.Range("J82").Value = 5003
$ws.Range("L82").Value = 5003
$ws.Range("N82").Value = -5725
$ws.Range("H83").Value = 65888.664
$ws.Range("I83").Value = 56000
$ws.Range("J83").Value = 78249.5
$ws.Range("K83").Value = 168000
$ws.Range("L83").Value = 234748.5
$ws.Range("M83").Value = -162384
$ws.Range("N83").Value = -245980.5
$ws.Range("H85").Value = 3500.4285
$ws.Range("J85").Value = 5003
$ws.Range("L85").Value = 5003
$ws.Range("N85").Value = -7499
$ws.Range("H86").Value = 59833.332
$ws.Range("J86").Value = 59833.332
$ws.Range("L86").Value = 59833.332
$ws.Range("N86").Value = -62205.332
$ws.Range("H89").Value = 59833.332
$ws.Range("J89").Value = 59833.332
$ws.Range("L89").Value = 179499.996
$ws.Range("N89").Value = -191355.996
$ws.Range("H100").Value = 3697.889
$ws.Range("I100").Value = 1950.4
$ws.Range("K100").Value = 1950.4
$ws.Range("M100").Value = -1409.4
$ws.Range("H113").Value = 4004.8948
$ws.Range("I113").Value = 3240.4707
$ws.Range("J113").Value = 10502.5
$ws.Range("K113").Value = 3240.4707
$ws.Range("L113").Value = 10502.5
$ws.Range("M113").Value = -1070.4707
$ws.Range("N113").Value = -14842.5
$ws.Range("H122").Value = 3754.0557
$ws.Range("I122").Value = 3382.4375
$ws.Range("K122").Value = 10147.3125
$ws.Range("M122").Value = -7697.3125
$ws.Range("H132").Value = 3530.2122
$ws.Range("I132").Value = 2930.3572
$ws.Range("J132").Value = 6889.4
$ws.Range("K132").Value = 8791.071599999999
$ws.Range("L132").Value = 20668.2
$ws.Range("M132").Value = -6261.071599999999
$ws.Range("N132").Value = -25728.2
$ws.Range("H136").Value = 4764.222
$ws.Range("I136").Value = 3334.5
$ws.Range("J136").Value = 7623.6665
$ws.Range("K136").Value = 10003.5
$ws.Range("L136").Value = 22870.9995
$ws.Range("M136").Value = -7453.5
$ws.Range("N136").Value = -27970.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20813.363
$ws.Range("J41").Value = 20813.363
$ws.Range("L41").Value = 20813.363
$ws.Range("N41").Value = -21593.363
$ws.Range("H62").Value = 79253.85000000001
$ws.Range("I62").Value = 204271.6
$ws.Range("J62").Value = 5714
$ws.Range("K62").Value = 204271.6
$ws.Range("L62").Value = 5714
$ws.Range("M62").Value = -203647.6
$ws.Range("N62").Value = -6962
$ws.Range("H65").Value = 79253.85000000001
$ws.Range("I65").Value = 204271.6
$ws.Range("J65").Value = 5714
$ws.Range("K65").Value = 1021358
$ws.Range("L65").Value = 28570
$ws.Range("M65").Value = -1018238
$ws.Range("N65").Value = -34810
$ws.Range("H70").Value = 29524.777
$ws.Range("J70").Value = 29465.375
$ws.Range("L70").Value = 29465.375
$ws.Range("N70").Value = -30095.375
$ws.Range("H73").Value = 29524.777
$ws.Range("J73").Value = 29465.375
$ws.Range("L73").Value = 29465.375
$ws.Range("N73").Value = -31649.375
$ws.Range("H81").Value = 1886.7333
$ws.Range("J81").Value = 4106.75
$ws.Range("L81").Value = 8213.5
$ws.Range("N81").Value = -10335.5
$ws.Range("H84").Value = 1886.7333
$ws.Range("J84").Value = 4106.75
$ws.Range("L84").Value = 41067.5
$ws.Range("N84").Value = -51675.5
$ws.Range("H96").Value = 3835.7334
$ws.Range("I96").Value = 4251.615
$ws.Range("J96").Value = 1132.5
$ws.Range("K96").Value = 4251.615
$ws.Range("L96").Value = 1132.5
$ws.Range("M96").Value = -2878.615
$ws.Range("N96").Value = -3878.5
$ws.Range("H100").Value = 1254.3704
$ws.Range("I100").Value = 1181.238
$ws.Range("K100").Value = 2362.476
$ws.Range("M100").Value = -1821.476
$ws.Range("H122").Value = 2572.6333
$ws.Range("I122").Value = 2588.111
$ws.Range("K122").Value = 7764.333
$ws.Range("M122").Value = -5314.333
$ws.Range("H132").Value = 10416.12
$ws.Range("I132").Value = 8601.450000000001
$ws.Range("J132").Value = 17674.8
$ws.Range("K132").Value = 25804.35
$ws.Range("L132").Value = 53024.39999999999
$ws.Range("M132").Value = -23274.35
$ws.Range("N132").Value = -58084.39999999999
$ws.Range("H133").Value = 79999.5
$ws.Range("J133").Value = 79999.5
$ws.Range("L133").Value = 79999.5
$ws.Range("N133").Value = -90119.5
$ws.Range("H136").Value = 5225.625
$ws.Range("I136").Value = 4766.7856
$ws.Range("K136").Value = 14300.3568
$ws.Range("M136").Value = -11750.3568

Write-Host "All updates applied."